# Auto update Excel log
# Appends new sensor-log rows (163-171) to the PIR, Humidity, and Temperature sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$rows = @(
  @(163, "'2026-01-28", "17:28:50", "17:00", "Bathroom", "No Motion", "Inactive"),
  @(164, "'2026-01-28", "17:28:55", "17:00", "Bathroom", "No Motion", "Inactive"),
  @(165, "'2026-01-28", "17:29:00", "17:00", "Bathroom", "No Motion", "Inactive"),
  @(166, "'2026-01-28", "17:29:05", "17:00", "Bathroom", "No Motion", "Inactive"),
  @(167, "'2026-01-28", "17:29:11", "17:00", "Bathroom", "No Motion", "Inactive"),
  @(168, "'2026-01-28", "17:29:15", "17:00", "Bathroom", "No Motion", "Inactive"),
  @(169, "'2026-01-28", "17:29:20", "17:00", "Bathroom", "No Motion", "Inactive"),
  @(170, "'2026-01-28", "17:29:25", "17:00", "Bathroom", "No Motion", "Inactive"),
  @(171, "'2026-01-28", "17:29:31", "17:00", "Bathroom", "No Motion", "Inactive")
)
foreach ($row in $rows) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
  $ws.Cells.Item($r, 6).Value = $row[6]
}

$ws = $wb.Worksheets.Item("Humidity")
$rows = @(
  @(163, "'2026-01-28", "17:28:45", "17:00", "Bathroom", "'86.9%", "Active"),
  @(164, "'2026-01-28", "17:28:57", "17:00", "Bathroom", "'87.8%", "Active"),
  @(165, "'2026-01-28", "17:29:01", "17:00", "Bathroom", "'87.8%", "Active"),
  @(166, "'2026-01-28", "17:29:06", "17:00", "Bathroom", "'86.9%", "Active"),
  @(167, "'2026-01-28", "17:29:09", "17:00", "Bathroom", "'87.8%", "Active"),
  @(168, "'2026-01-28", "17:29:17", "17:00", "Bathroom", "'86.9%", "Active"),
  @(169, "'2026-01-28", "17:29:21", "17:00", "Bathroom", "'87.7%", "Active"),
  @(170, "'2026-01-28", "17:29:26", "17:00", "Bathroom", "'86.8%", "Active"),
  @(171, "'2026-01-28", "17:29:29", "17:00", "Bathroom", "'87.7%", "Active")
)
foreach ($row in $rows) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
  $ws.Cells.Item($r, 6).Value = $row[6]
}

$ws = $wb.Worksheets.Item("Temperature")
$rows = @(
  @(163, "'2026-01-28", "17:28:46", "17:00", "Bathroom", "22.8C", "Active"),
  @(164, "'2026-01-28", "17:28:58", "17:00", "Bathroom", "22.8C", "Active"),
  @(165, "'2026-01-28", "17:29:02", "17:00", "Bathroom", "22.8C", "Active"),
  @(166, "'2026-01-28", "17:29:06", "17:00", "Bathroom", "22.8C", "Active"),
  @(167, "'2026-01-28", "17:29:10", "17:00", "Bathroom", "22.8C", "Active"),
  @(168, "'2026-01-28", "17:29:18", "17:00", "Bathroom", "22.8C", "Active"),
  @(169, "'2026-01-28", "17:29:22", "17:00", "Bathroom", "22.8C", "Active"),
  @(170, "'2026-01-28", "17:29:26", "17:00", "Bathroom", "22.8C", "Active"),
  @(171, "'2026-01-28", "17:29:30", "17:00", "Bathroom", "22.8C", "Active")
)
foreach ($row in $rows) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
  $ws.Cells.Item($r, 6).Value = $row[6]
}
